# Fremantle_stats.xlsx edit:
# The worksheet had its most-recent data column "JN" (unstyled, no explicit
# cell style) duplicated into three new columns (JO, JP, JQ) so that the
# single latest prediction is now repeated across JN:JQ, with JN, JO, JP
# taking on the explicit style that the rest of the sheet already uses and
# the new right-most column (JQ) inheriting the old "last column" unstyled
# look that JN used to have. No values actually change - it is a pure
# structural (copy/insert) edit that expands the sheet from A1:JN102 to
# A1:JQ102.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three copies of column JN immediately to its own right-hand side,
# each time re-copying the (still unstyled) original JN column. Inserting
# *at* JN pushes the existing JN content one column to the right each time,
# so after three passes the original JN values/format end up at JQ (still
# unstyled) while JN/JO/JP all carry the copied-cell style - exactly
# matching the target layout.
for ($i = 0; $i -lt 3; $i++) {
    $ws.Columns("JN").Copy()
    $ws.Columns("JN").Insert()
}
